$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price & 1h volume change), and a few row re-ranks

$ws.Range("D2").Value = "27.371.47"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "1.654.29"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.30"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.512"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.68"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0614"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0877"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "1.646.83"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.571"
$ws.Range("E14").Value = "  +3.91%  "
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.64"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "27.375.29"
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.63"
$ws.Range("E18").Value = "  -7.13%  "
$ws.Range("D19").Value = "0.0₃0728"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.38"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.35"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.73"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.17"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("D33").Value = "1.456.05"
$ws.Range("E33").Value = "  +3.06%  "
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.571"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.04"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.00"
$ws.Range("E43").Value = "  -6.18%  "
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.788"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.796.20"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.21"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0106"
$ws.Range("E49").Value = "  -3.24%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.76"
$ws.Range("E51").Value = "  +0.36%  "
